$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BOL": refresh the OrderId test values used by the reroute scripts
# ---------------------------------------------------------------------------
$bol = $wb.Worksheets.Item("BOL")

$bol.Range("A2").Value = "'51498462"
$bol.Range("A3").Value = "'51495312"
$bol.Range("A4").Value = "'51495312"

$bol.Range("D1").EntireColumn.ColumnWidth = 11.1667

$bol.Range("G5").Select()

# ---------------------------------------------------------------------------
# Sheet "Reroute Request": move the "Commercial" / "Amazon FBA Warehouse"
# sample rows further down the sheet and clear the stale OrderId / Way Bill #
# test values, then re-apply the boxed grid look to the whole table.
# ---------------------------------------------------------------------------
$rr = $wb.Worksheets.Item("Reroute Request")

# --- Row 9 (was row 3): Commercial test location -------------------------
$rr.Range("A9").Value = "Commercial"
$rr.Range("D9").Value = "Test Location 2"
$rr.Range("E9").Value = "#123, Moody"
$rr.Range("F9").Value = "MOODY"
$rr.Range("G9").Value = "AL"
$rr.Range("H9").Value = "'35004"
$rr.Range("I9").Value = "US"
$rr.Range("J9").Value = 0.5
$rr.Range("K9").Value = 0.66666666666666663
$rr.Range("J9").NumberFormat = "h:mm AM/PM"
$rr.Range("K9").NumberFormat = "h:mm AM/PM"

# --- Row 11 (was row 4): Amazon FBA Warehouse test location --------------
$rr.Range("A11").Value = "Amazon FBA Warehouse"
$rr.Range("D11").Value = "Test Location 3"
$rr.Range("E11").Value = "#321, Los Angeles"
$rr.Range("F11").Value = "LOS ANGELES"
$rr.Range("G11").Value = "CA"
$rr.Range("H11").Value = "'90001"
$rr.Range("I11").Value = "US"

# --- Clear out the old rows 3 & 4 (their content now lives in 9 / 11) ----
$rr.Range("A3:O3").Clear()
$rr.Range("A4:O4").Clear()

# --- Clear the OrderId / Way Bill # test data (kept blank going forward) -
$rr.Range("B2").ClearContents()
$rr.Range("C2").ClearContents()
$rr.Range("B9").ClearContents()
$rr.Range("C9").ClearContents()
$rr.Range("B11").ClearContents()
$rr.Range("C11").ClearContents()

# --- quotePrefix carries over automatically with the moved text values;
#     make sure the Way Bill # column keeps it for the moved rows too.
$rr.Range("C9").Value = "'"
$rr.Range("C9").ClearContents()
$rr.Range("C11").Value = "'"
$rr.Range("C11").ClearContents()

# --- Column widths ---------------------------------------------------------
$rr.Range("F1").EntireColumn.ColumnWidth = 11.6667
$rr.Range("O1").EntireColumn.ColumnWidth = 16.1667

# --- Give column O (orderReferenceID) the same boxed look on the data rows
$rr.Range("O2").Value = "'"
$rr.Range("O2").ClearContents()
$rr.Range("O9").Value = "'"
$rr.Range("O9").ClearContents()
$rr.Range("O11").Value = "'"
$rr.Range("O11").ClearContents()

# --- Re-draw the thin boxed border around the whole used range -----------
$rr.Range("A1:O11").Borders.LineStyle = 1

$rr.Range("I14").Select()
$rr.Activate()
